$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F ("dSF") based on repulled data / mean calculation
$updates = @{
    2  = -8
    3  = -7
    4  = -2
    6  = 5
    9  = -2
    11 = -12
    12 = -16
    13 = -6
    14 = -5
    16 = -4
    18 = -1
    20 = -1
    22 = -7
    29 = -6
    31 = 3
    34 = 3
    37 = 4
    38 = -1
    39 = 1
    41 = -1
    43 = -4
    44 = -6
    45 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
